$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A44").Value = "nu_leakSR"
$ws.Range("B44").Value = 0.2

$ws.Range("E44").Select()
